$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3: turn the literal FALSE boolean into a =FALSE() formula (keeps the same
# displayed value, matching what happened to C2's TRUE() earlier)
$ws.Range("C3").Formula = "=FALSE()"

# New row 4 data
$ws.Range("A4").Value = 186020
$ws.Range("C4").Value = $true
$ws.Range("C4").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("N4").Value = "1CxYvmK3CryGXBY9emG52Uiy2FggYASi8jrl9EUPbMgZLZkNqMLH1cpNZsYmEAqCadNXvEIL7MDTPhGrZQ2Xs4LprRQJZI27j4jxAtekCNxN17xY6l1akBTRsObEQFdK"

# Update the view: scroll so column B is the left-most visible column, and
# select F7 as the active cell
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("F7").Select()
